# Update cosinor_6_sine_01_.xlsx results with re-run CircaDB / CircadiPy analysis values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("E2").Value = [double]"25.8200000000006"
$ws.Range("H2").Value = [double]"5.586027796856133e-16"
$ws.Range("K2").Value = [double]"43.8677195619666"
$ws.Range("L2").Value = "[36.611000206572484, 51.124438917360706]"
$ws.Range("O2").Value = [double]"1.490605523324887"
$ws.Range("P2").Value = "[1.3145002294299628, 1.666710817219811]"
$ws.Range("S2").Value = [double]"58.99188405837011"
$ws.Range("T2").Value = "[54.37151505092708, 63.61225306581314]"
$ws.Range("W2").Value = [double]"19.69453453453499"
$ws.Range("X2").Value = [double]"18.97085085085128"
$ws.Range("Y2").Value = [double]"20.41821821821869"

# --- Row 3 ---
$ws.Range("B3").Value = [double]"0"
$ws.Range("E3").Value = [double]"22"
$ws.Range("G3").Value = [double]"4.088994698392412e-10"
$ws.Range("H3").Value = [double]"1.569155605036974e-09"
$ws.Range("I3").Value = [double]"0.2457135161537739"
$ws.Range("K3").Value = [double]"38.61458320221318"
$ws.Range("L3").Value = "[27.191955202473117, 50.037211201953234]"
$ws.Range("M3").Value = [double]"5.029670013811938e-10"
$ws.Range("N3").Value = [double]"5.029670013811938e-10"
$ws.Range("O3").Value = [double]"0.03144737390980801"
$ws.Range("P3").Value = "[-0.30818426431611545, 0.37107901213573147]"
$ws.Range("Q3").Value = [double]"0.8550223525430174"
$ws.Range("R3").Value = [double]"0.8550223525430174"
$ws.Range("S3").Value = [double]"58.87209449871644"
$ws.Range("T3").Value = "[51.43021425509957, 66.31397474233331]"
$ws.Range("W3").Value = [double]"21.88988988988989"
$ws.Range("X3").Value = [double]"20.7007007007007"
$ws.Range("Y3").Value = [double]"23.07907907907908"
